$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.605.91'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '3.319.76'
$ws.Range("E3").Value = '  +5.34%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.86'
$ws.Range("E5").Value = '  +3.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.06'
$ws.Range("E6").Value = '  +2.87%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.316.31'
$ws.Range("E8").Value = '  +5.33%  '
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("E11").Value = '  +3.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.71'
$ws.Range("E14").Value = '  +1.55%  '
$ws.Range("D15").Value = '3.867.01'
$ws.Range("E15").Value = '  +5.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.120'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '3.322.05'
$ws.Range("E17").Value = '  +5.45%  '
$ws.Range("D18").Value = '63.703.42'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.13'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("E22").Value = '  +4.79%  '
$ws.Range("E23").Value = '  +5.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.75'
$ws.Range("E24").Value = '  +5.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.01'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.25'
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.00'
$ws.Range("E32").Value = '  +8.14%  '
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.53'
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("E36").Value = '  +4.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.80'
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("D38").Value = '0.0₃0743'
$ws.Range("E38").Value = '  +6.58%  '
$ws.Range("E39").Value = '  +2.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '433.14'
$ws.Range("E40").Value = '  +4.12%  '
$ws.Range("D41").Value = '3.093.47'
$ws.Range("E41").Value = '  +5.64%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.34'
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("E44").Value = '  +4.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.263'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.73'
$ws.Range("E47").Value = '  +13.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.33'
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.998'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.01'
$ws.Range("E51").Value = '  +3.38%  '
